{"js": "// Auto-generated by gen/gen_js.py \u2014 do not hand-edit without regenerating.\n// Applies the resume bullet-text enhancements described in the commit diff.\n\nconst replacements = [\n  { index: 9, oldText: \"\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations\", newText: \"\u2022 Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\" },\n  { index: 10, oldText: \"\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics\", newText: \"\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis\" },\n  { index: 11, oldText: \"\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\", newText: \"\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\" },\n  { index: 12, oldText: \"\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\", newText: \"\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\" },\n  { index: 13, oldText: \"\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\", newText: \"\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\" },\n  { index: 14, oldText: \"\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\", newText: \"\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\" },\n  { index: 17, oldText: \"\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\", newText: \"\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\" },\n  { index: 18, oldText: \"\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\", newText: \"\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\" },\n  { index: 19, oldText: \"\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\", newText: \"\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\" },\n  { index: 20, oldText: \"\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products\", newText: \"\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products\" },\n  { index: 21, oldText: \"\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\", newText: \"\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\" },\n  { index: 33, oldText: \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research\", newText: \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions\" },\n  { index: 36, oldText: \"\u2022 Managed critical research operations for political campaigns\", newText: \"\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\" },\n  { index: 37, oldText: \"\u2022 Conducted comprehensive polling and demographic analysis\", newText: \"\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\" },\n  { index: 38, oldText: \"\u2022 Developed strategic recommendations based on data analysis\", newText: \"\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\" },\n  { index: 39, oldText: \"\u2022 Led research team in support of progressive political initiatives\", newText: \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs\" },\n  { index: 53, oldText: \"Political Research and Data Analysis\", newText: \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\" },\n  { index: 54, oldText: \"\u2022 Developed data analysis tools for political polling and research\", newText: \"\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\" },\n  { index: 55, oldText: \"\u2022 Built statistical models for voter behavior analysis\", newText: \"\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute\" },\n  { index: 56, oldText: \"\u2022 Created data visualization tools for research presentations\", newText: \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions\" },\n  { index: 57, oldText: \"\u2022 Supported senior researchers with technical analysis and reporting\", newText: \"\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\" },\n  { index: 59, oldText: \"Political Field Operations and Data Management\", newText: \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\" },\n  { index: 60, oldText: \"\u2022 Managed field operations for political campaigns and research projects\", newText: \"\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions\" },\n  { index: 61, oldText: \"\u2022 Developed data collection and management systems for field work\", newText: \"\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm\" },\n  { index: 62, oldText: \"\u2022 Trained field staff on data collection protocols and quality control\", newText: \"\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings\" },\n  { index: 63, oldText: \"\u2022 Analyzed field data to inform campaign strategy and research findings\", newText: \"\u2022 Created custom reports and data visualizations based on specific client requirements\" },\n];\n\nconst insertions = [\n  { anchorIndex: 39, anchorText: \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs\", newParagraphs: [\"\u2022 Managed comprehensive research operations for progressive political initiatives and candidates\"] },\n  { anchorIndex: 57, anchorText: \"\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\", newParagraphs: [\"\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps\", \"\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\"] },\n  { anchorIndex: 63, anchorText: \"\u2022 Created custom reports and data visualizations based on specific client requirements\", newParagraphs: [\"\u2022 Introduced mapping and geospatial analysis into standard reporting procedures\", \"\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\"] },\n];\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// --- Step 1: apply the 26 straight text replacements, by original\n// paragraph index. Each is verified against the expected \"old\" text\n// before being overwritten, so the script fails loudly instead of\n// silently touching the wrong paragraph if the document does not match\n// what we expect.\nfor (const r of replacements) {\n  const p = paras.items[r.index];\n  if (p.text !== r.oldText) {\n    throw new Error(\n      \"Paragraph \" + r.index + \" text mismatch.\\nExpected: \" + r.oldText +\n      \"\\nFound:    \" + p.text\n    );\n  }\n  p.insertText(r.newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Step 2: insert the 5 new bullet paragraphs. Insertions are\n// anchored on the paragraph (by original index) that must immediately\n// precede each new block; that paragraph's text is re-verified (now\n// post-replacement) before inserting after it. Because we insert in\n// ascending document order and always anchor on ORIGINAL indices (which\n// remain valid in `paras` \u2014 inserting paragraphs elsewhere does not\n// reshuffle the earlier, already-loaded collection array), no index\n// bookkeeping is required between insertions.\nfor (const ins of insertions) {\n  const anchor = paras.items[ins.anchorIndex];\n  anchor.load(\"text\");\n}\nawait context.sync();\n\nlet lastInserted = null;\nfor (const ins of insertions) {\n  const anchor = paras.items[ins.anchorIndex];\n  if (anchor.text !== ins.anchorText) {\n    throw new Error(\n      \"Anchor paragraph \" + ins.anchorIndex + \" text mismatch.\\nExpected: \" +\n      ins.anchorText + \"\\nFound:    \" + anchor.text\n    );\n  }\n  let insertAfter = anchor;\n  for (const newText of ins.newParagraphs) {\n    lastInserted = insertAfter.insertParagraph(newText, Word.InsertLocation.after);\n    insertAfter = lastInserted;\n  }\n}\nawait context.sync();\n", "ps1": "# Auto-generated by gen/gen_ps1.py \u2014 do not hand-edit without regenerating.\n# Applies the resume bullet-text enhancements described in the commit diff.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: 26 straight paragraph-text replacements, addressed by their\n# 1-based index in the ORIGINAL document (Word COM Paragraphs collection is\n# 1-based). Each paragraph's current text is verified against the expected\n# 'old' text before being overwritten, so the script fails loudly instead of\n# silently touching the wrong paragraph.\n$replacements = @(\n    @{ Index = 10; OldText = \"\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations\"; NewText = \"\u2022 Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\" }\n    @{ Index = 11; OldText = \"\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics\"; NewText = \"\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis\" }\n    @{ Index = 12; OldText = \"\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\"; NewText = \"\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\" }\n    @{ Index = 13; OldText = \"\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\"; NewText = \"\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\" }\n    @{ Index = 14; OldText = \"\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\"; NewText = \"\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\" }\n    @{ Index = 15; OldText = \"\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\"; NewText = \"\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\" }\n    @{ Index = 18; OldText = \"\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\"; NewText = \"\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\" }\n    @{ Index = 19; OldText = \"\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\"; NewText = \"\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\" }\n    @{ Index = 20; OldText = \"\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\"; NewText = \"\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\" }\n    @{ Index = 21; OldText = \"\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products\"; NewText = \"\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products\" }\n    @{ Index = 22; OldText = \"\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\"; NewText = \"\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\" }\n    @{ Index = 34; OldText = \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research\"; NewText = \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions\" }\n    @{ Index = 37; OldText = \"\u2022 Managed critical research operations for political campaigns\"; NewText = \"\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\" }\n    @{ Index = 38; OldText = \"\u2022 Conducted comprehensive polling and demographic analysis\"; NewText = \"\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\" }\n    @{ Index = 39; OldText = \"\u2022 Developed strategic recommendations based on data analysis\"; NewText = \"\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\" }\n    @{ Index = 40; OldText = \"\u2022 Led research team in support of progressive political initiatives\"; NewText = \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs\" }\n    @{ Index = 54; OldText = \"Political Research and Data Analysis\"; NewText = \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\" }\n    @{ Index = 55; OldText = \"\u2022 Developed data analysis tools for political polling and research\"; NewText = \"\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\" }\n    @{ Index = 56; OldText = \"\u2022 Built statistical models for voter behavior analysis\"; NewText = \"\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute\" }\n    @{ Index = 57; OldText = \"\u2022 Created data visualization tools for research presentations\"; NewText = \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions\" }\n    @{ Index = 58; OldText = \"\u2022 Supported senior researchers with technical analysis and reporting\"; NewText = \"\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\" }\n    @{ Index = 60; OldText = \"Political Field Operations and Data Management\"; NewText = \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\" }\n    @{ Index = 61; OldText = \"\u2022 Managed field operations for political campaigns and research projects\"; NewText = \"\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions\" }\n    @{ Index = 62; OldText = \"\u2022 Developed data collection and management systems for field work\"; NewText = \"\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm\" }\n    @{ Index = 63; OldText = \"\u2022 Trained field staff on data collection protocols and quality control\"; NewText = \"\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings\" }\n    @{ Index = 64; OldText = \"\u2022 Analyzed field data to inform campaign strategy and research findings\"; NewText = \"\u2022 Created custom reports and data visualizations based on specific client requirements\" }\n)\n\nforeach ($r in $replacements) {\n    $p = $d.Paragraphs.Item($r.Index)\n    $current = $p.Range.Text\n    # Range.Text includes the trailing paragraph mark (\\r); compare without it.\n    if ($current.Length -gt 0 -and $current.Substring($current.Length - 1) -eq [char]13) {\n        $current = $current.Substring(0, $current.Length - 1)\n    }\n    if ($current -ne $r.OldText) {\n        throw \"Paragraph $($r.Index) text mismatch.`nExpected: $($r.OldText)`nFound:    $current\"\n    }\n    $p.Range.Text = $r.NewText\n}\n\n# --- Step 2: insert the 5 new bullet paragraphs. Anchors are given by their\n# 1-based index in the ORIGINAL document; `$offset` accumulates how many\n# paragraphs earlier insertions in this loop have already added, so later\n# anchors keep resolving to the right paragraph even though the document\n# keeps growing as we go.\n$insertions = @(\n    @{ AnchorIndex = 40; AnchorText = \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs\"; NewParagraphs = @(\"\u2022 Managed comprehensive research operations for progressive political initiatives and candidates\") }\n    @{ AnchorIndex = 58; AnchorText = \"\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\"; NewParagraphs = @(\"\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps\", \"\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\") }\n    @{ AnchorIndex = 64; AnchorText = \"\u2022 Created custom reports and data visualizations based on specific client requirements\"; NewParagraphs = @(\"\u2022 Introduced mapping and geospatial analysis into standard reporting procedures\", \"\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\") }\n)\n\n$offset = 0\nforeach ($ins in $insertions) {\n    $anchorIndex = $ins.AnchorIndex + $offset\n    $anchor = $d.Paragraphs.Item($anchorIndex)\n    $current = $anchor.Range.Text\n    if ($current.Length -gt 0 -and $current.Substring($current.Length - 1) -eq [char]13) {\n        $current = $current.Substring(0, $current.Length - 1)\n    }\n    if ($current -ne $ins.AnchorText) {\n        throw \"Anchor paragraph $anchorIndex text mismatch.`nExpected: $($ins.AnchorText)`nFound:    $current\"\n    }\n    $insertAt = $anchorIndex\n    foreach ($newText in $ins.NewParagraphs) {\n        $d.Paragraphs.Item($insertAt).Range.InsertParagraphAfter()\n        $insertAt = $insertAt + 1\n        $d.Paragraphs.Item($insertAt).Range.Text = $newText\n    }\n    $offset = $offset + $ins.NewParagraphs.Count\n}\n\n"}
